# Applies the cryptos.xlsx price/volume(1h) update described in the commit
# "Updated cryptos list on Sun Dec 17 21:40:35 UTC 2023 with GitHub Actions".
# Row 12/13 additionally swap Polkadot <-> TRON (ranking order changed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.056.90'

$ws.Range("D3").Value = '2.238.71'
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''242.48'
$ws.Range("E5").Value = '  -1.15%  '

$ws.Range("D6").Value = '''0.616'
$ws.Range("E6").Value = '  -2.01%  '

$ws.Range("D7").Value = '''74.16'
$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.597'
$ws.Range("E9").Value = '  -3.81%  '

$ws.Range("D10").Value = '''42.11'
$ws.Range("E10").Value = '  -2.13%  '

$ws.Range("D11").Value = '''0.0951'
$ws.Range("E11").Value = '  -1.40%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '''6.90'
$ws.Range("E12").Value = '  -3.12%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.103'
$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("D14").Value = '2.571.90'
$ws.Range("E14").Value = '  +0.05%  '

$ws.Range("D15").Value = '''14.34'
$ws.Range("E15").Value = '  -0.94%  '

$ws.Range("E16").Value = '  -1.70%  '

$ws.Range("D17").Value = '2.257.74'
$ws.Range("E17").Value = '  +1.25%  '

$ws.Range("D18").Value = '41.933.02'
$ws.Range("E18").Value = '  -0.79%  '

$ws.Range("E19").Value = '  -5.82%  '

$ws.Range("E20").Value = '  +0.46%  '

$ws.Range("D21").Value = '''72.54'
$ws.Range("E21").Value = '  +0.57%  '

$ws.Range("D22").Value = '''11.11'
$ws.Range("E22").Value = '  +7.90%  '

$ws.Range("D23").Value = '''229.64'
$ws.Range("E23").Value = '  -0.70%  '

$ws.Range("D24").Value = '''2.03'
$ws.Range("E24").Value = '  -6.29%  '

$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").Value = '''11.32'
$ws.Range("E26").Value = '  -3.35%  '

$ws.Range("E27").Value = '  -0.86%  '

$ws.Range("D28").Value = '''2.27'
$ws.Range("E28").Value = '  -0.97%  '

$ws.Range("E29").Value = '  -0.85%  '

$ws.Range("D30").Value = '''167.50'
$ws.Range("E30").Value = '  +0.47%  '

$ws.Range("D31").Value = '''20.57'
$ws.Range("E31").Value = '  -1.78%  '

$ws.Range("E32").Value = '  -4.70%  '

$ws.Range("E33").Value = '  -0.81%  '

$ws.Range("D34").Value = '''30.24'
$ws.Range("E34").Value = '  +1.89%  '

$ws.Range("E35").Value = '  -0.70%  '

$ws.Range("D36").Value = '''0.110'
$ws.Range("E36").Value = '  -6.96%  '

$ws.Range("E37").Value = '  -4.09%  '

$ws.Range("E38").Value = '  -1.56%  '

$ws.Range("D39").Value = '''13.08'
$ws.Range("E39").Value = '  -1.09%  '

$ws.Range("D40").Value = '''2.12'
$ws.Range("E40").Value = '  -1.98%  '

$ws.Range("D41").Value = '''5.67'
$ws.Range("E41").Value = '  +0.24%  '

$ws.Range("D42").Value = '''64.50'
$ws.Range("E42").Value = '  +1.80%  '

$ws.Range("E43").Value = '  -1.89%  '

$ws.Range("D44").Value = '''8.70'
$ws.Range("E44").Value = '  -1.47%  '

$ws.Range("D45").Value = '''103.49'
$ws.Range("E45").Value = '  -1.93%  '

$ws.Range("E46").Value = '  -1.93%  '

$ws.Range("E47").Value = '  -0.27%  '

$ws.Range("E48").Value = '  -0.84%  '

$ws.Range("E49").Value = '  -2.19%  '

$ws.Range("E50").Value = '  -1.74%  '

$ws.Range("D51").Value = '2.447.16'
$ws.Range("E51").Value = '  +0.06%  '
